# Remove the trailing "Ver no Jupiter ..." / "© 2020 ..." footer block
# (and the blank paragraph that precedes it) that the site generator had
# appended to the course page. The blank paragraph that carries the
# page-break (directly before the footer block) and the blank paragraph
# that follows the footer block are left untouched.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph by its text (Range.Text carries
# a trailing paragraph-mark character, so match with StartsWith instead of
# an exact -eq comparison).
$verIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("Ver no Jupiter Salvar em pdf Salvar em docx")) {
        $verIndex = $i
        break
    }
}

if ($verIndex -eq -1) {
    throw "Could not find the 'Ver no Jupiter' paragraph"
}

# The blank paragraph immediately before it, and the copyright paragraph
# immediately after it, bound the block to delete.
$startPara = $d.Paragraphs.Item($verIndex - 1)
$endPara   = $d.Paragraphs.Item($verIndex + 1)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()
